# Rename all 50 worksheets in place (order/sheetId/r:id unchanged,
# only the sheet "name" attribute changes) to reflect a rerun of the
# underlying models (new random summary IDs baked into the sheet names).

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ26139037",
    "summ26278538",
    "summ26435099",
    "summ26594868",
    "summ26744524",
    "summ26893620",
    "summ27038407",
    "summ27288383",
    "summ27445156",
    "summ27597373",
    "summ27748098",
    "summ27933910",
    "summ28087653",
    "summ28237132",
    "summ28383969",
    "summ28542531",
    "summ28699757",
    "summ28895823",
    "summ29085247",
    "summ29274248",
    "summ29498642",
    "summ29691309",
    "summ29895021",
    "summ30093089",
    "summ30301864",
    "summ30520821",
    "summ30728923",
    "summ30932006",
    "summ31142078",
    "summ31352117",
    "summ31579884",
    "summ31793097",
    "summ32032519",
    "summ32234750",
    "summ32433813",
    "summ32618873",
    "summ32810108",
    "summ33006251",
    "summ33201721",
    "summ33391078",
    "summ33603506",
    "summ33809760",
    "summ34007825",
    "summ34209278",
    "summ34405738",
    "summ34614064",
    "summ34818059",
    "summ35023030",
    "summ35218145",
    "summ35412866"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $sheet = $wb.Worksheets.Item($i + 1)
    $sheet.Name = $newNames[$i]
}
